$p = $ppt.ActivePresentation
$master = $p.SlideMaster

# Update the date placeholder on the Slide Master
$master.Shapes.Item(3).TextFrame.TextRange.Text = "14.01.2013"

# Update the date placeholder on each Slide Layout
$layoutShapeIdx = @(3,3,3,4,6,2,1,4,4,3,3)
for ($i=1; $i -le $master.CustomLayouts.Count; $i++) {
    $lay = $master.CustomLayouts.Item($i)
    $idx = $layoutShapeIdx[$i-1]
    $lay.Shapes.Item($idx).TextFrame.TextRange.Text = "14.01.2013"
}

Write-Output "Done"
